$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")
$ws.Activate()

# Update the repaymentstrategy value (B17) from "RBI (India)" to the new scenario value
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Match the author's cursor/selection position after the edit
$ws.Range("B17").Select()
